# Add a "V-Sync Enabled" setting row to the strings sheet, right above
# "gamma" (the first graphics-related setting), shifting it and every
# later row down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "gamma" currently lives on row 25; insert a fresh row above it so
# "gamma" and everything below moves down to make room.
$ws.Rows.Item(25).Insert()

$ws.Range("A25").Value = "vsync enabled"
$ws.Range("B25").Value = "V-Sync Enabled"

# Match the row height used by the other simple (non-wrapping) setting rows.
$ws.Rows.Item(25).RowHeight = 13.4

# Leave the selection on the newly inserted row, like the authored commit.
$ws.Range("A26").Select() | Out-Null
